# Update data requirements input file.
# The "Forms" column (column D) in the TABLE worksheet is removed, shifting
# the remaining columns (Description / Data source, etc.) left by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TABLE")

# Delete the entire "Forms" column (D). This shifts columns E:G left to D:F
# and Excel automatically drops the now-unused shared strings (the old
# "Forms" header plus the per-row form codes such as "1RC", "3AR, 3CE, 3SU",
# etc.) when the file is saved.
$ws.Columns.Item(4).Delete()
